# Edit script: update "properties" sheet materials, add new rows, remove duplicates,
# sort alphabetically by name, freeze header row, and apply an AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

# --- 1. Rename materials that dropped the "-RES0" suffix ---
$ws.Cells.Item(14,1).Value2 = "Opaque Door panel_con-efficient"
$ws.Cells.Item(15,1).Value2 = "Overhead Door_con Panel-efficient"
$ws.Cells.Item(16,1).Value2 = "CP02 CARPET PAD-efficient"
$ws.Cells.Item(27,1).Value2 = "Std Opaque Door Panel-efficient"

$ws.Cells.Item(28,1).Value2 = "Opaque Door panel_con-non-standard"
$ws.Cells.Item(30,1).Value2 = "CP02 CARPET PAD-non-standard"
$ws.Cells.Item(41,1).Value2 = "Std Opaque Door Panel-non-standard"

$ws.Cells.Item(42,1).Value2 = "Opaque Door panel_con-standard"
$ws.Cells.Item(44,1).Value2 = "CP02 CARPET PAD-standard"
$ws.Cells.Item(55,1).Value2 = "Std Opaque Door Panel-standard"

$ws.Cells.Item(56,1).Value2 = "Opaque Door panel_con-ZEB"
$ws.Cells.Item(58,1).Value2 = "CP02 CARPET PAD-ZEB"
$ws.Cells.Item(69,1).Value2 = "Std Opaque Door Panel-ZEB"

# --- 2. Remove the per-variant "Air_Wall_Material-*-RES0" rows (bottom-up so row
#        indices of rows above stay valid) ---
$ws.Rows.Item(59).Delete()
$ws.Rows.Item(45).Delete()
$ws.Rows.Item(31).Delete()
$ws.Rows.Item(17).Delete()

# --- 3. Append the new consolidated rows at the bottom of the table ---
$lastRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($lastRow,1).Value2 = "CP02 CARPET PAD"
$ws.Cells.Item($lastRow,2).Value2 = 2500
$ws.Cells.Item($lastRow,3).Value2 = 0.006
$ws.Cells.Item($lastRow,4).Value2 = 0.9
$ws.Cells.Item($lastRow,5).Value2 = "USA"
$lastRow = $lastRow + 1

$ws.Cells.Item($lastRow,1).Value2 = "Air_Wall_Material"
$ws.Cells.Item($lastRow,2).Value2 = 2500
$ws.Cells.Item($lastRow,3).Value2 = 0.006
$ws.Cells.Item($lastRow,4).Value2 = 0.9
$ws.Cells.Item($lastRow,5).Value2 = "USA"

# --- 4. Sort the table (excluding header) alphabetically by column A ---
$fullRange = $ws.UsedRange
$fullRange.Sort($ws.Range("A1"), 1, $null, $null, 2, $null, 2, 1)

# --- 5. Freeze the header row, leaving the selection on the data rows below ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A3:XFD5").Select()

# --- 6. Apply an AutoFilter across the full table ---
$tableAddr = $ws.UsedRange.Address()
$ws.Range($tableAddr).AutoFilter()

# --- 7. Register the (hidden) _FilterDatabase defined name, scoped to the sheet ---
$fdbName = $ws.Names.Add("_xlnm._FilterDatabase", "=properties!" + $tableAddr)
$fdbName.Visible = $false
